$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row renames (E1:K1) -------------------------------------------
# "BusinessID" becomes the new unique identifier column "DynamiteInternalId",
# and the other Dynamite/Publishing field headers get their proper internal
# SharePoint field names.
$ws.Range("K1").Value = "DynamiteInternalId"
$ws.Range("E1").Value = "DynamiteNavigation"
$ws.Range("F1").Value = "PublishingPageContent"
$ws.Range("J1").Value = "DynamitePublishingStartDate"
$ws.Range("G1").Value = "DynamiteSummary"
$ws.Range("H1").Value = "PublishingPageImage"
$ws.Range("I1").Value = "DynamiteImageDescription"

# --- Column widths ----------------------------------------------------------
# Columns D:I used to share one uniform (default) width; split them out and
# size E:I (now holding the renamed Dynamite/Publishing columns) individually.
$ws.Columns.Item(5).ColumnWidth = 13.7366
$ws.Columns.Item(6).ColumnWidth = 21.3366
$ws.Columns.Item(7).ColumnWidth = 19.0
$ws.Columns.Item(8).ColumnWidth = 22.5
$ws.Columns.Item(9).ColumnWidth = 25.8366

# Column K (new unique id column) needs its own, wider, column so the new
# identifier values aren't truncated.
$ws.Columns.Item(11).ColumnWidth = 27.0

# --- Selection ---------------------------------------------------------------
$ws.Range("G6").Select()
